# Refresh the cryptos list: updated Price (D) and Volume(1h) (E) figures
# pulled from the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "43.015.72"
$ws.Range('E2').Value = "  -0.64%  "
$ws.Range('D3').Value = "2.342.10"
$ws.Range('E3').Value = "  +0.90%  "
$ws.Range('E4').Value = "  +0.01%  "
$ws.Range('D5').Value = "'302.82"
$ws.Range('E5').Value = "  +0.12%  "
$ws.Range('D6').Value = "'94.53"
$ws.Range('E6').Value = "  -3.34%  "
$ws.Range('E7').Value = "  -0.86%  "
$ws.Range('E8').Value = "  +0.00%  "
$ws.Range('E9').Value = "  -1.25%  "
$ws.Range('D10').Value = "'34.04"
$ws.Range('E10').Value = "  -4.13%  "
$ws.Range('D11').Value = "'0.0783"
$ws.Range('E11').Value = "  -1.44%  "
$ws.Range('D12').Value = "'18.64"
$ws.Range('E12').Value = "  -4.40%  "
$ws.Range('E13').Value = "  +1.63%  "
$ws.Range('D14').Value = "'6.74"
$ws.Range('E14').Value = "  -2.26%  "
$ws.Range('D15').Value = "2.705.76"
$ws.Range('E15').Value = "  +0.72%  "
$ws.Range('D16').Value = "2.364.76"
$ws.Range('E16').Value = "  +1.52%  "
$ws.Range('D17').Value = "'0.797"
$ws.Range('E17').Value = "  +0.81%  "
$ws.Range('D18').Value = "42.961.18"
$ws.Range('E18').Value = "  -0.57%  "
$ws.Range('D19').Value = "'12.07"
$ws.Range('E19').Value = "  -4.50%  "
$ws.Range('E20').Value = "  +2.14%  "
$ws.Range('D21').Value = "0.0₃0888"
$ws.Range('E21').Value = "  -1.03%  "
$ws.Range('D22').Value = "'67.94"
$ws.Range('E22').Value = "  -0.04%  "
$ws.Range('D23').Value = "'235.74"
$ws.Range('E23').Value = "  -0.65%  "
$ws.Range('D24').Value = "'2.22"
$ws.Range('E24').Value = "  -1.33%  "
$ws.Range('E25').Value = "  +0.02%  "
$ws.Range('E26').Value = "  -1.60%  "
$ws.Range('E27').Value = "  -1.59%  "
$ws.Range('D28').Value = "'2.35"
$ws.Range('E28').Value = "  +13.75%  "
$ws.Range('D29').Value = "'9.17"
$ws.Range('E29').Value = "  +0.60%  "
$ws.Range('D30').Value = "'31.46"
$ws.Range('E30').Value = "  -4.60%  "
$ws.Range('E31').Value = "  +0.02%  "
$ws.Range('D33').Value = "'0.0739"
$ws.Range('E33').Value = "  +5.15%  "
$ws.Range('D34').Value = "'17.23"
$ws.Range('E34').Value = "  -3.77%  "
$ws.Range('D37').Value = "'2.31"
$ws.Range('E37').Value = "  -1.64%  "
$ws.Range('D38').Value = "'0.100"
$ws.Range('E38').Value = "  -0.58%  "
$ws.Range('D39').Value = "'122.96"
$ws.Range('E39').Value = "  -25.13%  "
$ws.Range('E40').Value = "  -0.91%  "
$ws.Range('D41').Value = "'22.14"
$ws.Range('E41').Value = "  +15.35%  "
$ws.Range('E42').Value = "  -1.09%  "
$ws.Range('D43').Value = "1.939.13"
$ws.Range('E43').Value = "  -2.28%  "
$ws.Range('D44').Value = "'0.0281"
$ws.Range('E44').Value = "  +0.30%  "
$ws.Range('D45').Value = "'10.09"
$ws.Range('E45').Value = "  -5.36%  "
$ws.Range('D46').Value = "'2.09"
$ws.Range('E46').Value = "  +1.13%  "
$ws.Range('E47').Value = "  -2.82%  "
$ws.Range('D48').Value = "2.571.80"
$ws.Range('E48').Value = "  +0.77%  "
$ws.Range('E49').Value = "  +0.40%  "
$ws.Range('D50').Value = "'52.82"
$ws.Range('E50').Value = "  -1.99%  "
$ws.Range('D51').Value = "'71.68"
$ws.Range('E51').Value = "  -1.35%  "

# Rows 35 and 36 swapped rank order (ARBITRUM now ranks above RenderToken)
$ws.Range('B35').Value = "ARBITRUM"
$ws.Range('C35').Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D35').Value = "'1.83"
$ws.Range('E35').Value = "  +2.99%  "
$ws.Range('B36').Value = "RenderToken"
$ws.Range('C36').Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D36').Value = "'4.37"
$ws.Range('E36').Value = "  -2.61%  "
